# Auto-generated script applying odds updates for 2025-02-13 FlashScore workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 3.7
$ws.Range("M3").Value = 1.22
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 1.73
$ws.Range("Y3").Value = 2
$ws.Range("Z3").Value = 1.73
$ws.Range("G4").Value = 4.33
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2
$ws.Range("K4").Value = 1.83
$ws.Range("L4").Value = 2.88
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("S4").Value = 2.88
$ws.Range("T4").Value = 1.4
$ws.Range("U4").Value = 4.9
$ws.Range("V4").Value = 1.18
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 1.13
$ws.Range("Y4").Value = 1.67
$ws.Range("Z4").Value = 2.1
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 19
$ws.Range("AI4").Value = 5.5
$ws.Range("AK4").Value = 23
$ws.Range("AN4").Value = 8
$ws.Range("AO4").Value = 10
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 4.75
$ws.Range("J5").Value = 2.63
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88
$ws.Range("S5").Value = 2.6
$ws.Range("T5").Value = 1.48
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 1.14
$ws.Range("AD5").Value = 7.5
$ws.Range("AF5").Value = 15
$ws.Range("AM5").Value = 9
$ws.Range("AN5").Value = 21
$ws.Range("AQ5").Value = 41
$ws.Range("H6").Value = 2.7
$ws.Range("I6").Value = 3.25
$ws.Range("L6").Value = 4.33
$ws.Range("M6").Value = 1.17
$ws.Range("N6").Value = 4.75
$ws.Range("O6").Value = 1.73
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = 2.41
$ws.Range("R6").Value = 1.58
$ws.Range("U6").Value = 5.8
$ws.Range("V6").Value = 1.14
$ws.Range("Y6").Value = 1.75
$ws.Range("Z6").Value = 2.05
$ws.Range("AA6").Value = 2.5
$ws.Range("AB6").Value = 1.5
$ws.Range("AF6").Value = 29
$ws.Range("AI6").Value = 4.75
$ws.Range("AJ6").Value = 6
$ws.Range("AK6").Value = 23
$ws.Range("AQ6").Value = 41
$ws.Range("G8").Value = 1.95
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 4.33
$ws.Range("J8").Value = 2.75
$ws.Range("L8").Value = 4.75
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("Y8").Value = 1.53
$ws.Range("Z8").Value = 2.38
$ws.Range("AF8").Value = 17
$ws.Range("AG8").Value = 19
$ws.Range("AK8").Value = 17
$ws.Range("AP8").Value = 41
$ws.Range("T9").Value = 1.44
$ws.Range("U9").Value = 4.1
$ws.Range("V9").Value = 1.22
$ws.Range("G15").Value = 4.2
$ws.Range("J15").Value = 4.3
$ws.Range("L15").Value = 2.15
$ws.Range("S15").Value = 1.62
$ws.Range("T15").Value = 2.02
$ws.Range("W15").Value = 2.45
$ws.Range("X15").Value = 1.42
$ws.Range("AC15").Value = 14
$ws.Range("AD15").Value = 25
$ws.Range("AF15").Value = 65
$ws.Range("AJ15").Value = 7.5
$ws.Range("AK15").Value = 14
$ws.Range("AM15").Value = 8.5
$ws.Range("AN15").Value = 9.25
$ws.Range("AP15").Value = 14
